# Weekly fruit/vegetable price update:
# Insert two new rows of fresh weekly data at the top of the
# "Terminal La Palmera de La Serena - Pepino dulce" price block
# (rows 456-457), pushing the existing rows down by two (they keep
# all their original values, just shifted to rows 458-484).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block (A456:R457). This
# shifts the previous rows 456..482 down to 458..484 and carries
# their content/formatting with them - no need to touch those rows.
$ws.Range("A456:R457").Insert()

# --- New row 456 : Primera ---
$ws.Range("A456").Value = 8
$ws.Range("B456").Value = "Terminal La Palmera de La Serena"
$ws.Range("C456").Value = "Coquimbo"
$ws.Range("D456").Value = 44783
$ws.Range("E456").Value = 4
$ws.Range("F456").Value = 100112043
$ws.Range("G456").Value = "Pepino dulce"
$ws.Range("H456").Value = "Cultivar IV Región"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 400
$ws.Range("K456").Value = 13000
$ws.Range("L456").Value = 13500
$ws.Range("M456").Value = 13250
$ws.Range("N456").Value = "`$/bandeja 18 kilos"
$ws.Range("O456").Value = "Provincia de Limarí"
$ws.Range("P456").Value = 736
$ws.Range("Q456").Value = 18
$ws.Range("R456").Value = "Hortaliza"

# --- New row 457 : Segunda ---
$ws.Range("A457").Value = 8
$ws.Range("B457").Value = "Terminal La Palmera de La Serena"
$ws.Range("C457").Value = "Coquimbo"
$ws.Range("D457").Value = 44783
$ws.Range("E457").Value = 4
$ws.Range("F457").Value = 100112043
$ws.Range("G457").Value = "Pepino dulce"
$ws.Range("H457").Value = "Cultivar IV Región"
$ws.Range("I457").Value = "Segunda"
$ws.Range("J457").Value = 200
$ws.Range("K457").Value = 11000
$ws.Range("L457").Value = 11500
$ws.Range("M457").Value = 11250
$ws.Range("N457").Value = "`$/bandeja 18 kilos"
$ws.Range("O457").Value = "Provincia de Limarí"
$ws.Range("P457").Value = 625
$ws.Range("Q457").Value = 18
$ws.Range("R457").Value = "Hortaliza"
